$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps storing values as text (matches the
# original inline-string cells) instead of Excel auto-converting
# numeric-looking strings like "1.003" into real numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.485.23"
$ws.Range("E2").Value = "  -0.60%  "
$ws.Range("D3").Value = "1.827.21"
$ws.Range("E3").Value = "  -1.71%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -1.06%  "
$ws.Range("D5").Value = "333.30"
$ws.Range("E5").Value = "  -0.29%  "
$ws.Range("D6").Value = "1.004"
$ws.Range("D7").Value = "0.4574"
$ws.Range("E7").Value = "  -0.79%  "
$ws.Range("D8").Value = "0.3834"
$ws.Range("E8").Value = "  -1.78%  "
$ws.Range("D9").Value = "46.07"
$ws.Range("E9").Value = "  -0.92%  "
$ws.Range("D10").Value = "0.07853"
$ws.Range("E10").Value = "  -1.09%  "
$ws.Range("D11").Value = "0.9589"
$ws.Range("E11").Value = "  -4.01%  "
$ws.Range("E12").Value = "  -2.57%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.823.89"
$ws.Range("E13").Value = "  -2.45%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "5.838"
$ws.Range("E14").Value = "  -1.79%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "7.051"
$ws.Range("E15").Value = "  -2.15%  "
$ws.Range("D16").Value = "1.005"
$ws.Range("E16").Value = "  -0.92%  "
$ws.Range("D17").Value = "89.55"
$ws.Range("E17").Value = "  +1.52%  "
$ws.Range("D18").Value = "0.06589"
$ws.Range("E18").Value = "  -1.98%  "
$ws.Range("E19").Value = "  -1.84%  "
$ws.Range("D20").Value = "17.11"
$ws.Range("E20").Value = "  -0.46%  "
$ws.Range("E21").Value = "  -0.85%  "
$ws.Range("D22").Value = "27.475.65"
$ws.Range("E22").Value = "  -0.67%  "
$ws.Range("D23").Value = "5.289"
$ws.Range("E23").Value = "  -2.59%  "
$ws.Range("D24").Value = "10.80"
$ws.Range("E24").Value = "  -1.31%  "
$ws.Range("E25").Value = "  -2.18%  "
$ws.Range("D26").Value = "159.24"
$ws.Range("E26").Value = "  -0.14%  "
$ws.Range("D27").Value = "2.041.55"
$ws.Range("E27").Value = "  -2.13%  "
$ws.Range("E28").Value = "  -1.38%  "
$ws.Range("D29").Value = "2.045"
$ws.Range("E29").Value = "  -3.85%  "
$ws.Range("D30").Value = "5.275"
$ws.Range("E30").Value = "  -3.31%  "
$ws.Range("D31").Value = "117.87"
$ws.Range("E31").Value = "  -3.19%  "
$ws.Range("D32").Value = "0.09370"
$ws.Range("E32").Value = "  -0.15%  "
$ws.Range("D33").Value = "0.9296"
$ws.Range("E33").Value = "  -4.45%  "
$ws.Range("D34").Value = "3.573"
$ws.Range("E34").Value = "  -1.49%  "
$ws.Range("D35").Value = "5.218"
$ws.Range("E35").Value = "  -1.54%  "
$ws.Range("E36").Value = "  -0.75%  "
$ws.Range("D37").Value = "0.05902"
$ws.Range("E37").Value = "  -1.66%  "
$ws.Range("D38").Value = "0.02180"
$ws.Range("E38").Value = "  -2.43%  "
$ws.Range("D39").Value = "8.095"
$ws.Range("E39").Value = "  -3.12%  "
$ws.Range("E40").Value = "  -0.83%  "
$ws.Range("D41").Value = "1.146"
$ws.Range("E41").Value = "  -3.89%  "
$ws.Range("D42").Value = "0.5734"
$ws.Range("E42").Value = "  -3.33%  "
$ws.Range("D43").Value = "0.1819"
$ws.Range("E43").Value = "  -2.57%  "
$ws.Range("D44").Value = "9.909"
$ws.Range("E44").Value = "  -3.99%  "
$ws.Range("D45").Value = "1.265"
$ws.Range("E45").Value = "  +1.51%  "
$ws.Range("D46").Value = "0.5406"
$ws.Range("E46").Value = "  -3.37%  "
$ws.Range("D47").Value = "11.77"
$ws.Range("E47").Value = "  -3.26%  "
$ws.Range("D48").Value = "1.884"
$ws.Range("E48").Value = "  -1.48%  "
$ws.Range("D49").Value = "0.06964"
$ws.Range("E49").Value = "  +3.80%  "
$ws.Range("D50").Value = "110.34"
$ws.Range("E50").Value = "  -1.16%  "
$ws.Range("E51").Value = "  -33.10%  "